# Updated GSC export files: append two more days of data (2025-11-06 and
# 2025-11-07) to the "Chart" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @{ Row = 32; Date = "2025-11-06"; NonHttps = 0; Https = 100 },
    @{ Row = 33; Date = "2025-11-07"; NonHttps = 0; Https = 94 }
)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    # Force a leading apostrophe so the engine stores this as literal text
    # (a shared string) instead of auto-converting the "yyyy-MM-dd" looking
    # text into a date serial number. ClearFormats() then drops the
    # quote-prefix cell style that the apostrophe leaves behind, so the new
    # cell ends up using the same default/general style as its neighbours.
    $dateCell.Value = "'" + $r.Date
    $dateCell.ClearFormats()

    $ws.Cells.Item($r.Row, 2).Value = $r.NonHttps
    $ws.Cells.Item($r.Row, 3).Value = $r.Https
}

Write-Output "Appended 2025-11-06 and 2025-11-07 rows to Chart sheet"
